# Weekly update: a new "Apio" (celery) price observation (week of
# date-serial 44509) is inserted ahead of the existing history, pushing
# the prior rows down by two. Two rows are added - one for "Primera"
# quality and one for "Segunda" quality - mirroring the layout of the
# surrounding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows right before the current row 414.
$ws.Rows("414:415").Insert()

# New row 414 - "Primera" quality.
$ws.Cells.Item(414, 1).Value = 6
$ws.Cells.Item(414, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(414, 3).Value = "Metropolitana"
$ws.Cells.Item(414, 4).Value = 44509
$ws.Cells.Item(414, 5).Value = 13
$ws.Cells.Item(414, 6).Value = 100112017
$ws.Cells.Item(414, 7).Value = "Apio"
$ws.Cells.Item(414, 8).Value = "Americana (o)"
$ws.Cells.Item(414, 9).Value = "Primera"
$ws.Cells.Item(414, 10).Value = 2100
$ws.Cells.Item(414, 11).Value = 6000
$ws.Cells.Item(414, 12).Value = 7000
$ws.Cells.Item(414, 13).Value = 6571
$ws.Cells.Item(414, 14).Value = "`$/docena de matas"
$ws.Cells.Item(414, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(414, 16).Value = 1095
$ws.Cells.Item(414, 17).Value = 6
$ws.Cells.Item(414, 18).Value = "Hortaliza"

# New row 415 - "Segunda" quality.
$ws.Cells.Item(415, 1).Value = 6
$ws.Cells.Item(415, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(415, 3).Value = "Metropolitana"
$ws.Cells.Item(415, 4).Value = 44509
$ws.Cells.Item(415, 5).Value = 13
$ws.Cells.Item(415, 6).Value = 100112017
$ws.Cells.Item(415, 7).Value = "Apio"
$ws.Cells.Item(415, 8).Value = "Americana (o)"
$ws.Cells.Item(415, 9).Value = "Segunda"
$ws.Cells.Item(415, 10).Value = 600
$ws.Cells.Item(415, 11).Value = 5000
$ws.Cells.Item(415, 12).Value = 5000
$ws.Cells.Item(415, 13).Value = 5000
$ws.Cells.Item(415, 14).Value = "`$/docena de matas"
$ws.Cells.Item(415, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(415, 16).Value = 833
$ws.Cells.Item(415, 17).Value = 6
$ws.Cells.Item(415, 18).Value = "Hortaliza"
